$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A, Q, R between row 2 and row 4
$ws.Range("A2").Value = 111525236
$ws.Range("Q2").Value = 404471.921986955
$ws.Range("R2").Value = 6706721.507764764

$ws.Range("A4").Value = 111525237
$ws.Range("Q4").Value = 404472.9612160316
$ws.Range("R4").Value = 6706723.452812355
